$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# Row 11: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A11:AB11").Clear()
$ws.Range("A11").Value = " "
$ws.Range("B11").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"PANELS-(PANELS)`", `"Right side panel -(Right side panel )`", `"EXTERNAL PANEL RUSTED CORROSION`", `"`", `"`", `"`"]"

# Row 12: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A12:AB12").Clear()
$ws.Range("A12").Value = " "
$ws.Range("B12").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"FLOOR BOARD DIRTY BY MILK POWDER DUST & ODOUR`", `"`", `"`", `"`"]"

# Row 14: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A14:AB14").Clear()
$ws.Range("A14").Value = " "
$ws.Range("B14").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"PANELS-(PANELS)`", `"Right side panel -(Right side panel )`", `"EXTERNAL PANEL RUSTED CORROSION BADLY.`", `"`", `"`", `"`"]"

# Row 15: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A15:AB15").Clear()
$ws.Range("A15").Value = " "
$ws.Range("B15").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"FLOOR BOARD-(FLOOR BOARD)`", `"FLOOR BOARD DIRTY BY MILK POWDER DUST & SCRATCHED`", `"`", `"`", `"`"]"

# Row 22: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A22:AB22").Clear()
$ws.Range("A22").Value = " "
$ws.Range("B22").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"F/B DIRTY BY SAND DUST & ODOUR `", `"`", `"`", `"`"]"

# Row 24: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A24:AB24").Clear()
$ws.Range("A24").Value = " "
$ws.Range("B24").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"F/B DIRTY BY SCRATCHED,OIL STAIN,BLACK STAIN & DUST.`", `"`", `"`", `"`"]"

# Row 27: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A27:AB27").Clear()
$ws.Range("A27").Value = " "
$ws.Range("B27").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"F/B DIRTY BY DUST .`", `"`", `"`", `"`"]"

# Row 33: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A33:AB33").Clear()
$ws.Range("A33").Value = " "
$ws.Range("B33").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"PANELS-(PANELS)`", `"Right side panel -(Right side panel )`", `"INTERNAL PANEL DIRTY & RUSTED CORROSION`", `"`", `"`", `"`"]"

# Row 34: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A34:AB34").Clear()
$ws.Range("A34").Value = " "
$ws.Range("B34").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"PANELS-(PANELS)`", `"Right side panel -(Right side panel )`", `"EXTERNAL PANEL RUSTED CORROSION. BADLY`", `"`", `"`", `"`"]"

# Row 35: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A35:AB35").Clear()
$ws.Range("A35").Value = " "
$ws.Range("B35").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"FLOOR BOARD-(FLOOR BOARD)`", `"FLOOR BOARD DIRTY BY MILK POWDER DUST & ODOUR`", `"`", `"`", `"`"]"

# Row 38: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A38:AB38").Clear()
$ws.Range("A38").Value = " "
$ws.Range("B38").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"PANELS-(PANELS)`", `"Right side panel -(Right side panel )`", `"INTERNAL PANEL INK DIRTY .`", `"`", `"`", `"`"]"

# Row 39: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A39:AB39").Clear()
$ws.Range("A39").Value = " "
$ws.Range("B39").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"F/B BADLY SAND DUST & ODOUR.`", `"`", `"`", `"`"]"

# Row 41: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A41:AB41").Clear()
$ws.Range("A41").Value = " "
$ws.Range("B41").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"F/B DIRTY BY TYER MARK,BLACK STAIN & DUST.`", `"`", `"`", `"`"]"

# Row 43: replace the old "duplicate" damage-text-only row with the
# collapsed [A=" ", B="[nil,...]"] representation.
$ws.Range("A43:AB43").Clear()
$ws.Range("A43").Value = " "
$ws.Range("B43").Value = "[nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, nil, `"FLOORS-(F)`", `"Threshold plate-(Threshold plate)`", `"F/BOARD DIRTY BY DUST.`", `"`", `"`", `"`"]"

# Column width tweaks (hidden helper column B grows; column X shrinks).
$ws.Columns.Item(2).ColumnWidth = 226.71428571428572
$ws.Columns.Item(24).ColumnWidth = 46.42857142857143

